{"js": "// Replace the old multi-run \"Dates de la campanya 2018 en qu\u00e8 usem la\n// constel\u00b7laci\u00f3 Perseus 30 d'octubre al novembre 8 i 29 de novembre de\n// desembre 8\" paragraphs with a single plain run reading\n// \"Dates de la campanya Taurus: 16-25 de gener\" (no run formatting),\n// across every occurrence in the document body.\n\nconst OLD_SNIPPET = \"Dates de la campanya\";\nconst NEW_TEXT = \"Dates de la campanya Taurus: 16-25 de gener\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  // Match the specific \"campaign dates\" paragraphs (old text mentions the\n  // constellation Perseus / 2018 dates); ignore unrelated prose that merely\n  // contains the phrase \"dates de la campanya\".\n  if (text.indexOf(OLD_SNIPPET) !== -1 && text.indexOf(\"constel\") !== -1) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\nfor (const paragraph of targets) {\n  const range = paragraph.getRange();\n  // Clear every run (and its formatting) from the paragraph, leaving an\n  // empty paragraph, then insert the new text as a brand-new, unformatted\n  // run \u2014 matching the target OOXML which has a bare <w:r><w:t>.\n  range.clear();\n  range.insertText(NEW_TEXT, \"Start\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the old multi-run \"Dates de la campanya 2018 en que usem la\n# constel\u00b7lacio Perseus 30 d'octubre al novembre 8 i 29 de novembre de\n# desembre 8\" paragraphs with a single plain run reading\n# \"Dates de la campanya Taurus: 16-25 de gener\" (no run formatting),\n# across every occurrence in the document body.\n\n$d = $word.ActiveDocument\n\n$NEW_TEXT = \"Dates de la campanya Taurus: 16-25 de gener\"\n\n# Collect the 1-based paragraph indexes of every \"campaign dates\" paragraph\n# first (old text mentions the constellation/\"constel...\" so this ignores\n# unrelated prose that merely contains the phrase \"dates de la campanya\").\n$targets = @()\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.IndexOf(\"Dates de la campanya\") -ge 0 -and $t.IndexOf(\"constel\") -ge 0) {\n        $targets += $i\n    }\n    $i = $i + 1\n}\n\nforeach ($idx in $targets) {\n    $p = $d.Paragraphs.Item($idx)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark so the delete doesn't merge this\n    # paragraph with the next one; this empties the paragraph's content\n    # while keeping the <w:p> (and its pPr/sectPr) intact.\n    $body = $d.Range($r.Start, $r.End - 1)\n    $body.Delete()\n    # Insert the replacement as a brand-new, unformatted run (no rPr) \u2014\n    # matching the target OOXML's bare <w:r><w:t>.\n    $body.InsertAfter($NEW_TEXT)\n}\n"}
